$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.8
$ws.Range("B3").Value = "'FALSE"
$ws.Range("C5").Value = 0.9
$ws.Range("C7").Value = 1
$ws.Range("B8").Value = "'TRUE"
$ws.Range("C10").Value = 1
$ws.Range("C11").Value = 0.8
